# Auto-generated Excel COM-interop script to apply market-price data refresh
# across the Ixion_Profits leve-profit workbook (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 396.78125
$ws.Range("I15").Value = 396.78125
$ws.Range("K15").Value = 1190.34375
$ws.Range("M15").Value = -1021.34375
$ws.Range("H64").Value = 3914.75
$ws.Range("I64").Value = 4209.2593
$ws.Range("J64").Value = 3303.077
$ws.Range("K64").Value = 4209.2593
$ws.Range("L64").Value = 3303.077
$ws.Range("M64").Value = -3961.2593
$ws.Range("N64").Value = -3799.077
$ws.Range("H67").Value = 3914.75
$ws.Range("I67").Value = 4209.2593
$ws.Range("J67").Value = 3303.077
$ws.Range("K67").Value = 4209.2593
$ws.Range("L67").Value = 3303.077
$ws.Range("M67").Value = -3351.2593
$ws.Range("N67").Value = -5019.077
$ws.Range("H132").Value = 1706.2941
$ws.Range("I132").Value = 1756.6875
$ws.Range("K132").Value = 5270.0625
$ws.Range("M132").Value = -2740.0625
$ws.Range("H138").Value = 2063.9854
$ws.Range("I138").Value = 1462.1818
$ws.Range("J138").Value = 2351.8044
$ws.Range("K138").Value = 4386.5454
$ws.Range("L138").Value = 7055.4132
$ws.Range("M138").Value = 753.4546
$ws.Range("N138").Value = -17335.4132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8662.107
$ws.Range("I32").Value = 7027.45
$ws.Range("J32").Value = 28278
$ws.Range("K32").Value = 7027.45
$ws.Range("L32").Value = 28278
$ws.Range("M32").Value = -6740.45
$ws.Range("N32").Value = -28852
$ws.Range("H134").Value = 53207.25
$ws.Range("J134").Value = 53207.25
$ws.Range("L134").Value = 53207.25
$ws.Range("N134").Value = -63347.25
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 15869
$ws.Range("I105").Value = 25235.445
$ws.Range("J105").Value = 3826.4285
$ws.Range("K105").Value = 25235.445
$ws.Range("L105").Value = 3826.4285
$ws.Range("M105").Value = -23488.445
$ws.Range("N105").Value = -7320.4285
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3847989.2
$ws.Range("I16").Value = 5918551.5
$ws.Range("J16").Value = 2659
$ws.Range("K16").Value = 5918551.5
$ws.Range("L16").Value = 2659
$ws.Range("M16").Value = -5918264.5
$ws.Range("N16").Value = -3233
$ws.Range("H31").Value = 9237.362999999999
$ws.Range("J31").Value = 9237.362999999999
$ws.Range("L31").Value = 9237.362999999999
$ws.Range("N31").Value = -9827.362999999999
$ws.Range("H34").Value = 9237.362999999999
$ws.Range("J34").Value = 9237.362999999999
$ws.Range("L34").Value = 9237.362999999999
$ws.Range("N34").Value = -9641.362999999999
$ws.Range("H50").Value = 19333.334
$ws.Range("J50").Value = 19333.334
$ws.Range("L50").Value = 19333.334
$ws.Range("N50").Value = -20583.334
$ws.Range("H54").Value = 8000
$ws.Range("J54").Value = 8000
$ws.Range("L54").Value = 8000
$ws.Range("N54").Value = -9316
$ws.Range("H62").Value = 5060.3076
$ws.Range("I62").Value = 5460.5
$ws.Range("J62").Value = 4420
$ws.Range("K62").Value = 5460.5
$ws.Range("L62").Value = 4420
$ws.Range("M62").Value = -4836.5
$ws.Range("N62").Value = -5668
$ws.Range("H65").Value = 5060.3076
$ws.Range("I65").Value = 5460.5
$ws.Range("J65").Value = 4420
$ws.Range("K65").Value = 27302.5
$ws.Range("L65").Value = 22100
$ws.Range("M65").Value = -24182.5
$ws.Range("N65").Value = -28340
$ws.Range("H113").Value = 3847989.2
$ws.Range("I113").Value = 5918551.5
$ws.Range("J113").Value = 2659
$ws.Range("K113").Value = 5918551.5
$ws.Range("L113").Value = 2659
$ws.Range("M113").Value = -5916381.5
$ws.Range("N113").Value = -6999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 3846276.8
$ws.Range("J12").Value = 96.8
$ws.Range("L12").Value = 290.4
$ws.Range("N12").Value = -636.4
$ws.Range("H98").Value = 12500757
$ws.Range("I98").Value = 350
$ws.Range("J98").Value = 20001000
$ws.Range("K98").Value = 1050
$ws.Range("L98").Value = 60003000
$ws.Range("M98").Value = 448
$ws.Range("N98").Value = -60005996
$ws.Range("H131").Value = 1755414.9
$ws.Range("J131").Value = 1069.0638
$ws.Range("L131").Value = 3207.1914
$ws.Range("N131").Value = -13287.1914
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 41794.332
$ws.Range("I69").Value = 32182
$ws.Range("J69").Value = 46600.5
$ws.Range("K69").Value = 32182
$ws.Range("L69").Value = 46600.5
$ws.Range("M69").Value = -31433
$ws.Range("N69").Value = -48098.5
$ws.Range("H72").Value = 41794.332
$ws.Range("I72").Value = 32182
$ws.Range("J72").Value = 46600.5
$ws.Range("K72").Value = 96546
$ws.Range("L72").Value = 139801.5
$ws.Range("M72").Value = -92802
$ws.Range("N72").Value = -147289.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 62501944
$ws.Range("I68").Value = 2123.625
$ws.Range("J68").Value = 125001760
$ws.Range("K68").Value = 2123.625
$ws.Range("L68").Value = 125001760
$ws.Range("M68").Value = -1374.625
$ws.Range("N68").Value = -125003258
$ws.Range("H71").Value = 62501944
$ws.Range("I71").Value = 2123.625
$ws.Range("J71").Value = 125001760
$ws.Range("K71").Value = 10618.125
$ws.Range("L71").Value = 625008800
$ws.Range("M71").Value = -6874.125
$ws.Range("N71").Value = -625016288
$ws.Range("H132").Value = 27789642
$ws.Range("I132").Value = 83362584
$ws.Range("J132").Value = 3169.125
$ws.Range("K132").Value = 250087752
$ws.Range("L132").Value = 9507.375
$ws.Range("M132").Value = -250085222
$ws.Range("N132").Value = -14567.375
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7493
$ws.Range("I62").Value = 7607.846
$ws.Range("J62").Value = 6000
$ws.Range("K62").Value = 7607.846
$ws.Range("L62").Value = 6000
$ws.Range("M62").Value = -6983.846
$ws.Range("N62").Value = -7248
$ws.Range("H65").Value = 7493
$ws.Range("I65").Value = 7607.846
$ws.Range("J65").Value = 6000
$ws.Range("K65").Value = 38039.23
$ws.Range("L65").Value = 30000
$ws.Range("M65").Value = -34919.23
$ws.Range("N65").Value = -36240
$ws.Range("H81").Value = 3100
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 3100
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 6200
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -8322
$ws.Range("H84").Value = 3100
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 3100
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 31000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -41608
$ws.Range("H132").Value = 1486.7234
$ws.Range("I132").Value = 980.4865
$ws.Range("J132").Value = 3359.8
$ws.Range("K132").Value = 2941.4595
$ws.Range("L132").Value = 10079.4
$ws.Range("M132").Value = -411.4594999999999
$ws.Range("N132").Value = -15139.4
$ws.Range("H133").Value = 40070
$ws.Range("J133").Value = 40070
$ws.Range("L133").Value = 40070
$ws.Range("N133").Value = -50190
$ws.Range("H135").Value = 40915
$ws.Range("J135").Value = 40915
$ws.Range("L135").Value = 40915
$ws.Range("N135").Value = -51055
$ws.Range("H138").Value = 44164.5
$ws.Range("J138").Value = 44164.5
$ws.Range("L138").Value = 44164.5
$ws.Range("N138").Value = -54444.5
$ws.Range("H141").Value = 57254.375
$ws.Range("J141").Value = 57254.375
$ws.Range("L141").Value = 57254.375
$ws.Range("N141").Value = -67614.375
